$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "79_1023": rerun of analysis - start_type Entry -> Notification/
# diagnosis, sanatorium Yes -> No, and a new c1b (K) column of zeros.
# ---------------------------------------------------------------------
$ws1023 = $wb.Worksheets.Item("79_1023")

for ($r = 2; $r -le 29; $r++) {
    $ws1023.Cells.Item($r, 5).Value = "Notification/diagnosis"
    $ws1023.Cells.Item($r, 6).Value = "No"
    $ws1023.Cells.Item($r, 11).Value = 0
}

$ws1023.Columns.Item(5).ColumnWidth = 19.0
$ws1023.Columns.Item(6).ColumnWidth = 9.15

# ---------------------------------------------------------------------
# Sheet "79_1023_sev": same rerun, plus a new all-cause-mortality-removed
# P column of zeros (no longer derived from all-cause mortality).
# ---------------------------------------------------------------------
$wsSev = $wb.Worksheets.Item("79_1023_sev")

for ($r = 2; $r -le 4; $r++) {
    $wsSev.Cells.Item($r, 5).Value = "Notification/diagnosis"
    $wsSev.Cells.Item($r, 6).Value = "No"
    $wsSev.Cells.Item($r, 11).Value = 0
    $wsSev.Cells.Item($r, 16).Value = 0
}

$wsSev.Columns.Item(5).ColumnWidth = 19.0

# ---------------------------------------------------------------------
# UI selection state updates. "79_1023" must end up the active tab
# (activeTab moves from 15 -> 14), so select/activate it last.
# ---------------------------------------------------------------------
$ws75 = $wb.Worksheets.Item("75")
$ws75.Range("E2").Select() | Out-Null

$wsSev.Range("E2").Select() | Out-Null

$ws1023.Range("E2:E29").Select() | Out-Null
